# Apply the cell updates described by the commit diff.
# A leading apostrophe forces Excel to treat the assigned value as
# literal text (avoiding numeric auto-conversion of values such as
# "523.37" or "0.999"); re-applying the "Normal" style afterwards
# clears the transient quote-prefix formatting Excel applies, so
# the cell style matches the original (unstyled) cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.654.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.10%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.105.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.31%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'523.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.00%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'140.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.03%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.01%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.100.40"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.16%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.07%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'7.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.30%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.57%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +2.43%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.637.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.15%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +1.50%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'26.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.63%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0000164"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.53%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'57.709.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.10%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.101.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.01%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.51%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'12.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.01%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'8.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.81%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'336.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.54%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.20%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +2.71%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +1.19%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.09%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.12%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0₃0922"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.44%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +3.47%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.04%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.93%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +2.43%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +1.05%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'20.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.10%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'155.93"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.70%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'4.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.55%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'6.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.46%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'26.87"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.10%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +1.36%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -1.38%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +12.89%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'RenzoRestakedETH"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'3.147.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.27%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Mantle"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.687"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +4.68%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'3.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.01%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'36.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.54%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.09%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.299.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.73%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.24%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.974"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +5.82%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'20.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.04%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +2.21%  "
$ws.Range("E51").Style = "Normal"
